$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells we touch to remain Text so values like
# "1.00", "258.50", or "0.0000103" are not re-interpreted as numbers.
$dCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.751.15'
$ws.Range("E2").Value = '  +3.29%  '

$ws.Range("D3").Value = '2.187.22'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '258.50'
$ws.Range("E5").Value = '  +2.08%  '

$ws.Range("D6").Value = '81.32'
$ws.Range("E6").Value = '  +8.20%  '

$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  +1.68%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = '0.595'
$ws.Range("E9").Value = '  +2.09%  '

$ws.Range("D10").Value = '43.15'
$ws.Range("E10").Value = '  +5.03%  '

$ws.Range("D11").Value = '0.0919'
$ws.Range("E11").Value = '  +0.80%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.97'
$ws.Range("E12").Value = '  +2.97%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  +2.05%  '

$ws.Range("D14").Value = '2.515.53'
$ws.Range("E14").Value = '  +0.39%  '

$ws.Range("D15").Value = '14.29'
$ws.Range("E15").Value = '  +0.91%  '

$ws.Range("D16").Value = '2.175.14'
$ws.Range("E16").Value = '  -0.40%  '

$ws.Range("D17").Value = '0.773'
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").Value = '43.619.60'
$ws.Range("E18").Value = '  +3.19%  '

$ws.Range("D19").Value = '0.0000103'
$ws.Range("E19").Value = '  +0.48%  '

$ws.Range("D20").Value = '70.37'
$ws.Range("E20").Value = '  -0.42%  '

$ws.Range("D21").Value = '5.92'
$ws.Range("E21").Value = '  +0.97%  '

$ws.Range("D22").Value = '2.38'
$ws.Range("E22").Value = '  +9.49%  '

$ws.Range("D23").Value = '230.34'
$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("D24").Value = '8.96'
$ws.Range("E24").Value = '  -6.40%  '

$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("B26").Value = 'InjectiveProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D26").Value = '41.10'
$ws.Range("E26").Value = '  +10.83%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '10.63'
$ws.Range("E27").Value = '  +1.67%  '

$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("D31").Value = '172.56'
$ws.Range("E31").Value = '  +1.91%  '

$ws.Range("D32").Value = '20.35'
$ws.Range("E32").Value = '  +1.78%  '

$ws.Range("D33").Value = '0.0866'
$ws.Range("E33").Value = '  +6.25%  '

$ws.Range("D34").Value = '5.26'
$ws.Range("E34").Value = '  +2.92%  '

$ws.Range("D35").Value = '0.115'
$ws.Range("E35").Value = '  +7.08%  '

$ws.Range("D36").Value = '0.122'
$ws.Range("E36").Value = '  +1.35%  '

$ws.Range("D37").Value = '4.48'
$ws.Range("E37").Value = '  +5.51%  '

$ws.Range("D38").Value = '0.0355'
$ws.Range("E38").Value = '  +6.28%  '

$ws.Range("D39").Value = '13.29'
$ws.Range("E39").Value = '  +12.09%  '

$ws.Range("D40").Value = '2.87'
$ws.Range("E40").Value = '  +17.92%  '

$ws.Range("D41").Value = '2.09'
$ws.Range("E41").Value = '  +1.69%  '

$ws.Range("D42").Value = '62.35'
$ws.Range("E42").Value = '  +4.70%  '

$ws.Range("D43").Value = '5.46'
$ws.Range("E43").Value = '  +5.72%  '

$ws.Range("D44").Value = '0.199'
$ws.Range("E44").Value = '  +0.94%  '

$ws.Range("D45").Value = '101.24'
$ws.Range("E45").Value = '  -1.72%  '

$ws.Range("D46").Value = '0.0978'
$ws.Range("E46").Value = '  +0.70%  '

$ws.Range("D47").Value = '8.22'
$ws.Range("E47").Value = '  -0.57%  '

$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").Value = '1.17'
$ws.Range("E48").Value = '  +3.47%  '

$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '1.11'
$ws.Range("E49").Value = '  +1.86%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '1.53'
$ws.Range("E50").Value = '  +27.63%  '

$ws.Range("D51").Value = '0.436'
$ws.Range("E51").Value = '  -7.43%  '

# Restore default (Normal) style on the Price cells so only the value changed
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}